$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.374.65"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "2.095.85"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.667"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.55"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "62.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +7.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "2.399.59"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.835"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "2.094.28"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.53%  "
$ws.Range("D18").Value = "37.274.43"
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("E20").Value = "  +14.39%  "
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  +4.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +6.40%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("E29").Value = "  +4.44%  "
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  +25.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.21%  "
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0624"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0904"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +154.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.01%  "
$ws.Range("E43").Value = "  +6.87%  "
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0975"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "99.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "1.333.33"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("E50").Value = "  +7.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.59%  "
